$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (B1:S1): a column was inserted, several labels were
# renamed/reordered, and two new trailing columns ("sulfoxide","cycle")
# were added before the final "weight" column.
$headers = @(
    "('C', 'O', 1)",
    "('C', 'O', 1.5)",
    "('N', 'H', 1)",
    "('S', 'O', 2)",
    "('C', 'S', 1)",
    "('N', 'C', 1)",
    "('C', 'C', 1)",
    "('C', 'C', 1.5)",
    "('C', 'H', 1)",
    "C",
    "H",
    "N",
    "O",
    "amines",
    "ether",
    "sulfoxide",
    "cycle",
    "weight"
)

# New data rows (B2:S6), replacing the previous (now stale) cluster data.
$data = @(
    @(0.07000000000000001, 0.063, 3.197, 0.082, 0.136, 9.842000000000001, 9.042999999999999, 1.604, 26.709, 15.932, 29.992, 4.546, 0.968, 3.468, 0.061, 0.082, 1.817, 307.586),
    @(0.608, 0.145, 2.184, 0.191, 0.575, 3.038, 4.949, 0.358, 9.709, 7.772, 12.247, 2.579, 1.555, 1.22, 0.08699999999999999, 0.191, 0.971, 209.485),
    @(2.77, 0.06, 0.749, 0.08, 0.11, 4.526, 24.051, 0.426, 40.298, 25.528, 41.361, 1.538, 3.423, 1.537, 1.18, 0.08, 4.398, 409.29),
    @(1.036, 0.415, 1.12, 0.13, 0.499, 0.628, 5.29, 18.261, 13.075, 22.594, 14.658, 2.274, 1.954, 0.212, 0.335, 0.13, 1.124, 407.514),
    @(0.982, 0.372, 1.748, 0.07099999999999999, 0.117, 4.011, 8.304, 23.589, 31.967, 32.151, 34.039, 2.589, 0.672, 1.446, 0.55, 0.07099999999999999, 0.976, 490.443)
)

$numCols = $headers.Length
$numRows = 1 + $data.Length
$arr = New-Object 'object[,]' $numRows,$numCols

for ($c = 0; $c -lt $numCols; $c++) {
    $arr[0, $c] = $headers[$c]
}
for ($r = 0; $r -lt $data.Length; $r++) {
    for ($c = 0; $c -lt $numCols; $c++) {
        $arr[$r + 1, $c] = $data[$r][$c]
    }
}

# Write all values (B1:S6) in one shot via a properly-shaped 2D array.
$ws.Range("B1:S6").Value = $arr

# The previous used range only went out to column R; give the new S column
# (both header and data) the same header style ("s=1", bold/centered/boxed)
# that's already applied across row 1 and column A, by copying the format
# from the neighbouring R1 cell that already carries it.
$ws.Range("R1").Copy()
$ws.Range("S1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
